$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the four names next to their labels (previously-blank value cells) ---
$ws.Range("F9").Value  = "Mario Vidal"          # Supervisor turno día:
$ws.Range("F10").Value = "Marcelino Vasquez"    # Supervisores turno noche:
$ws.Range("F11").Value = "Jaime Urra"           # Ito Planta Turno día:
$ws.Range("F12").Value = "Victor Cortes"        # Ito Planta Turno noche:

# --- Convert the ISO datetime string to a localized date/time string ---
# Old "Hora Inicio:" / "Hora Termino:" value cells (F14/F15) referenced the
# ISO timestamp; they become blank.
$ws.Range("F14").ClearContents()
$ws.Range("F15").ClearContents()

# The "Precalentamiento" row's start/end timestamps (D20/E20) get reformatted
# from ISO 8601 to dd/mm/yyyy hh:mm:ss text.
$ws.Range("D20").Value = "16/04/2017 00:32:13"
$ws.Range("E20").Value = "16/04/2017 00:32:13"

# --- New "Precalentamiento2" row (21) gains the same start/end timestamp and
#     a 100% progress value, mirroring row 20 ---
$ws.Range("D21").Value = "16/04/2017 00:32:13"
$ws.Range("E21").Value = "16/04/2017 00:32:13"

# "100%" must land as literal text (matching F20's cell exactly), not get
# auto-converted to a percentage number by input parsing. Instead of typing
# the text (which would make Excel reinterpret "100%" as 1 formatted as a
# percentage), copy F20 verbatim - it already holds the literal text "100%" -
# onto F21 so both the value and the style/format match F20 exactly.
$ws.Range("F20").Copy()
$ws.Range("F21").PasteSpecial(-4104)
$excel.CutCopyMode = 0
